# chartink_screener.xlsx - "break out stock.yaml completed"
#
# 1) On sheet "10per change", rows 11-19 had their bsecode (column D)
#    stored as text; they become real numbers.
# 2) A fresh block of rows (20-28) is appended, duplicating rows 11-19's
#    data (bsecode stored as text again, as it originally was) with an
#    updated "Date Time" stamp (17:18:12 instead of 17:10:12).
# 3) The same pattern repeats on sheet "DND 3 V 0.3": row 3's bsecode
#    becomes numeric, and a new row 4 duplicates it with text bsecode and
#    the later timestamp.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Cell, $Text)
    # Force a literal/text cell (not Excel's "looks like a number -> number"
    # auto-conversion) the same way a user would in the UI: a leading
    # apostrophe. Then drop back to the "Normal" style so the cell doesn't
    # keep the quote-prefix formatting flag hanging off it.
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "10per change"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("10per change")

$rows = @(
    @{ A = 1; B = "LT";         C = "Larsen & Toubro Limited";                D = "500510"; E = 0.17;  F = 3409;    G = 10372458 },
    @{ A = 2; B = "LODHA";      C = "Macrotech Developers Ltd";               D = "543287"; E = -0.53; F = 1296.85; G = 1797849 },
    @{ A = 3; B = "ATGL";       C = "Adani Total Gas Ltd";                    D = "542066"; E = 3.03;  F = 936.25;  G = 6670432 },
    @{ A = 4; B = "SBIN";       C = "State Bank Of India";                    D = "500112"; E = 1.88;  F = 789.75;  G = 74256082 },
    @{ A = 5; B = "ADANIPOWER"; C = "Adani Power Limited";                    D = "533096"; E = 0.51;  F = 726.65;  G = 34537620 },
    @{ A = 6; B = "CGPOWER";    C = "CG Power and Industrial Solutions Ltd";  D = "500093"; E = 0.1;   F = 627.65;  G = 10674892 },
    @{ A = 7; B = "PAYTM";      C = "One 97 Communications Ltd";              D = "543396"; E = -4.91; F = 339.85;  G = 7404922 },
    @{ A = 8; B = "POWERGRID";  C = "Power Grid Corporation Of India Limited";D = "532898"; E = 0.96;  F = 298.8;   G = 45312613 },
    @{ A = 9; B = "GAIL";       C = "Gail (india) Limited";                   D = "532155"; E = 2.55;  F = 195.15;  G = 49797002 }
)

# 1) Rows 11-19: bsecode (column D) switches from text to a real number.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 11 + $i
    $ws1.Cells.Item($r, 4).Value = [double]$rows[$i].D
}

# 2) Append rows 20-28: same data, bsecode back to text, new timestamp.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 20 + $i
    $row = $rows[$i]
    $ws1.Cells.Item($r, 1).Value = $row.A
    $ws1.Cells.Item($r, 2).Value = $row.B
    $ws1.Cells.Item($r, 3).Value = $row.C
    Set-TextValue $ws1.Cells.Item($r, 4) $row.D
    $ws1.Cells.Item($r, 5).Value = $row.E
    $ws1.Cells.Item($r, 6).Value = $row.F
    $ws1.Cells.Item($r, 7).Value = $row.G
    $ws1.Cells.Item($r, 8).Value = "05/06/2024 17:18:12"
}

# ---------------------------------------------------------------------
# Sheet "DND 3 V 0.3"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DND 3 V 0.3")

# 1) Row 3: bsecode switches from text to a real number.
$ws3.Cells.Item(3, 4).Value = 500331

# 2) Append row 4: same data, bsecode back to text, new timestamp.
$ws3.Cells.Item(4, 1).Value = 1
$ws3.Cells.Item(4, 2).Value = "PIDILITIND"
$ws3.Cells.Item(4, 3).Value = "Pidilite Industries Limited"
Set-TextValue $ws3.Cells.Item(4, 4) "500331"
$ws3.Cells.Item(4, 5).Value = 3.4
$ws3.Cells.Item(4, 6).Value = 3166.2
$ws3.Cells.Item(4, 7).Value = 632880
$ws3.Cells.Item(4, 8).Value = "05/06/2024 17:18:12"
